# Junction_Flooding_305.xlsx — "custom accuracy" edit:
#   1. Row 5's numeric sensor readings are rounded down to 2-decimal
#      ("custom") precision (were stored with 3 decimals).
#   2. Row 6 (the next timestamped reading) is dropped, shrinking the
#      sheet's used range from A1:AH6 to A1:AH5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New (2-decimal) values for every data column of row 5, keyed by cell ref.
$rounded = [ordered]@{
    "B5"  = 18.74
    "C5"  = 13.74
    "D5"  = 1.18
    "E5"  = 40.45
    "F5"  = 33.48
    "G5"  = 14.75
    "H5"  = 58.03
    "I5"  = 22.69
    "J5"  = 10.04
    "K5"  = 15.05
    "L5"  = 16.28
    "M5"  = 17.08
    "N5"  = 4.71
    "O5"  = 14.66
    "P5"  = 20.85
    "Q5"  = 12.37
    "R5"  = 0.86
    "S5"  = 0.78
    "T5"  = 215.87
    "U5"  = 41.03
    "V5"  = 13.53
    "W5"  = 27.54
    "X5"  = 14.64
    "Y5"  = 1.84
    "Z5"  = 27.99
    "AA5" = 11.95
    "AB5" = 10.63
    "AC5" = 12.49
    "AD5" = 17.07
    "AE5" = 0.56
    "AF5" = 52.58
    "AG5" = 7.62
    "AH5" = 16.92
}

foreach ($ref in $rounded.Keys) {
    $ws.Range($ref).Value = $rounded[$ref]
}

# Drop row 6 entirely, shifting nothing else below it up (there is nothing
# below it) and shrinking the sheet's dimension accordingly.
$ws.Range("A6").EntireRow.Delete()
